$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string "Down" used in Y2 - set via cell value
$ws.Range("X2").Value = -0.059997999999993112
$ws.Range("Y2").Value = "Down"

# New row 3 data
$ws.Range("A3").Value = 42648.663923611108
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 38
$ws.Range("E3").Value = 5611
$ws.Range("F3").Value = 1008
$ws.Range("G3").Value = 63
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 83
$ws.Range("J3").Value = 16
$ws.Range("K3").Value = 10457
$ws.Range("L3").Value = 137
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 46
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = "Noun"
$ws.Range("Q3").Value = 48.098617091043238
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.096699999999999994
$ws.Range("S3").NumberFormat = $ws.Range("S2").NumberFormat
$ws.Range("T3").Value = 0.027400000000000001
$ws.Range("T3").NumberFormat = $ws.Range("T2").NumberFormat
$ws.Range("U3").Value = 4.8
$ws.Range("V3").Value = 2.2799999999999998
$ws.Range("W3").Value = 0
